$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("for json")

# Add the new "dsmavg"-style delta column E: E{n} = C{n} - C{n-1}, for rows 10..164.
# Setting .Formula on the whole range lets relative references auto-adjust per row
# and groups them into a shared formula, matching how Excel fills a column down.
$ws.Range("E10:E164").Formula = "=C10-C9"

# Move/restore the active selection to D10, as recorded after the edit.
$ws.Range("D10").Select() | Out-Null
